$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Functional Requirements")
$ws2 = $wb.Worksheets.Item("Non-Functional Requirements")

# Non-Functional Requirements sheet: the uptime/availability requirement
# (row 8) is promoted from "Want" to "Must", and its summary/description
# text is reworded accordingly.
$ws2.Range("E8").Value = "The TVM must be available for use at least 99% of the time."
$ws2.Range("C8").Value = "The TVM must have an uptime at least 99%"
$ws2.Range("D8").Value = "Must"

# Functional Requirements sheet: the last use case (row 17) used to be the
# "report issue" feature; it is reworked into a "feedback section" use case.
$ws1.Range("C17").Value = "The TVM should have a feedback section on the menu"
$ws1.Range("E17").Value = "The TVM will have a feedback section for user to send their experience on the TVM."

# Update the remembered selection on each affected sheet.
$ws1.Range("E17").Select()
$ws2.Range("D12").Select()
